$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 19

$ws.Cells.Item($row, 1).Value = 42601.898865740739
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($row, 2).Value = "Bag"
$ws.Cells.Item($row, 3).Value = 8666
$ws.Cells.Item($row, 4).Value = 13447
$ws.Cells.Item($row, 5).Value = 1638
$ws.Cells.Item($row, 6).Value = 180
$ws.Cells.Item($row, 7).Value = 87
$ws.Cells.Item($row, 8).Value = 66
$ws.Cells.Item($row, 9).Value = 32
$ws.Cells.Item($row, 10).Value = 2
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 100
$ws.Cells.Item($row, 13).Value = 0
